# Update odds values on Sheet1 to match the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 3 (Santa Fe vs Bucaramanga)
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5

# Row 4 (Atl. Nacional vs Deportes Tolima)
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.25
$ws.Range("S4").Value = 1.95
$ws.Range("T4").Value = 1.9
$ws.Range("W4").Value = 3.4
$ws.Range("X4").Value = 1.3
$ws.Range("AA4").Value = 1.91
$ws.Range("AB4").Value = 1.8
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 13
$ws.Range("AK4").Value = 19
$ws.Range("AN4").Value = 13
$ws.Range("AO4").Value = 26
$ws.Range("AP4").Value = 17
